# Add data for 2023-04-21
# Updates column J (year 2023 year-to-date totals) across all affected sheets.
$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Range("J2").Value = 2036
$ws.Range("J3").Value = 2119
$ws.Range("J4").Value = 485
$ws.Range("J5").Value = 155
$ws.Range("J6").Value = 2654
$ws.Range("J7").Value = 7449

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Range("J6").Value = 79
$ws.Range("J8").Value = 466
$ws.Range("J11").Value = 101
$ws.Range("J15").Value = 96
$ws.Range("J18").Value = 86
$ws.Range("J19").Value = 247
$ws.Range("J20").Value = 153
$ws.Range("J23").Value = 69
$ws.Range("J24").Value = 26
$ws.Range("J27").Value = 44
$ws.Range("J29").Value = 418
$ws.Range("J33").Value = 312
$ws.Range("J34").Value = 43
$ws.Range("J36").Value = 112
$ws.Range("J37").Value = 249
$ws.Range("J43").Value = 75
$ws.Range("J46").Value = 25
$ws.Range("J48").Value = 68
$ws.Range("J51").Value = 100
$ws.Range("J52").Value = 175
$ws.Range("J54").Value = 153
$ws.Range("J55").Value = 88
$ws.Range("J57").Value = 40
$ws.Range("J60").Value = 47
$ws.Range("J61").Value = 14
$ws.Range("J63").Value = 35
$ws.Range("J65").Value = 194
$ws.Range("J67").Value = 273
$ws.Range("J70").Value = 14
$ws.Range("J75").Value = 28
$ws.Range("J76").Value = 111
$ws.Range("J77").Value = 54
$ws.Range("J79").Value = 222
$ws.Range("J83").Value = 178
$ws.Range("J84").Value = 75
$ws.Range("J85").Value = 349
$ws.Range("J89").Value = 81
$ws.Range("J90").Value = 83
$ws.Range("J91").Value = 89
$ws.Range("J95").Value = 113
$ws.Range("J98").Value = 46
$ws.Range("J99").Value = 96
$ws.Range("J100").Value = 14
$ws.Range("J101").Value = 7449

# Sheet 3: South Shore
$ws = $wb.Worksheets.Item(3)
$ws.Range("J3").Value = 141
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 349

# Sheet 5: Little Village
$ws = $wb.Worksheets.Item(5)
$ws.Range("J2").Value = 40
$ws.Range("J3").Value = 53
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 175

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item(6)
$ws.Range("J2").Value = 37
$ws.Range("J3").Value = 20
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 101

# Sheet 7: Austin
$ws = $wb.Worksheets.Item(7)
$ws.Range("J2").Value = 146
$ws.Range("J3").Value = 154
$ws.Range("J6").Value = 132
$ws.Range("J7").Value = 466

# Sheet 10: Uptown
$ws = $wb.Worksheets.Item(10)
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 81

# Sheet 14: Grand Crossing
$ws = $wb.Worksheets.Item(14)
$ws.Range("J2").Value = 71
$ws.Range("J3").Value = 89
$ws.Range("J5").Value = 10
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 249

# Sheet 15: Woodlawn
$ws = $wb.Worksheets.Item(15)
$ws.Range("J2").Value = 31
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 96

# Sheet 16: North Lawndale
$ws = $wb.Worksheets.Item(16)
$ws.Range("J3").Value = 107
$ws.Range("J7").Value = 273

# Sheet 18: South Deering
$ws = $wb.Worksheets.Item(18)
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 75

# Sheet 19: New City
$ws = $wb.Worksheets.Item(19)
$ws.Range("J2").Value = 57
$ws.Range("J3").Value = 52
$ws.Range("J7").Value = 194

# Sheet 20: South Chicago
$ws = $wb.Worksheets.Item(20)
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 178

# Sheet 21: West Pullman
$ws = $wb.Worksheets.Item(21)
$ws.Range("J2").Value = 40
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 113

# Sheet 22: Garfield Park
$ws = $wb.Worksheets.Item(22)
$ws.Range("J2").Value = 79
$ws.Range("J3").Value = 93
$ws.Range("J6").Value = 111
$ws.Range("J7").Value = 312

# Sheet 24: Loop
$ws = $wb.Worksheets.Item(24)
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 153

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Range("J2").Value = 126
$ws.Range("J3").Value = 146
$ws.Range("J6").Value = 109
$ws.Range("J7").Value = 418

# Sheet 26: Chatham
$ws = $wb.Worksheets.Item(26)
$ws.Range("J3").Value = 68
$ws.Range("J7").Value = 247

# Sheet 28: Lake View
$ws = $wb.Worksheets.Item(28)
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 68

# Sheet 29: River North
$ws = $wb.Worksheets.Item(29)
$ws.Range("J6").Value = 63
$ws.Range("J7").Value = 111

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item(30)
$ws.Range("J2").Value = 22
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 79

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item(32)
$ws.Range("J4").Value = 16
$ws.Range("J6").Value = 140

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item(36)
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 88

# Sheet 37: Dunning
$ws = $wb.Worksheets.Item(37)
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 26

# Sheet 38: Jefferson Park
$ws = $wb.Worksheets.Item(38)
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 25

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item(39)
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 69

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item(40)
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 89

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item(42)
$ws.Range("J2").Value = 63
$ws.Range("J3").Value = 83
$ws.Range("J7").Value = 222

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item(44)
$ws.Range("J4").Value = 18
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 153

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item(45)
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 86

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item(47)
$ws.Range("J2").Value = 38
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 112

# Sheet 49: Wrigleyville
$ws = $wb.Worksheets.Item(49)
$ws.Range("J2").Value = 4
$ws.Range("J6").Value = 14

# Sheet 50: Garfield Ridge
$ws = $wb.Worksheets.Item(50)
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 43

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item(53)
$ws.Range("J4").Value = 2

# Sheet 54: Brighton Park
$ws = $wb.Worksheets.Item(54)
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 96

# Sheet 55: Wicker Park
$ws = $wb.Worksheets.Item(55)
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 46

# Sheet 67: O'Hare
$ws = $wb.Worksheets.Item(67)
$ws.Range("J3").Value = 3
$ws.Range("J7").Value = 14

# Sheet 71: Edgewater
$ws = $wb.Worksheets.Item(71)
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 44

# Sheet 73: Pullman
$ws = $wb.Worksheets.Item(73)
$ws.Range("J3").Value = 6
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 28

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item(74)
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 83

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item(75)
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 100

# Sheet 77: Mckinley Park
$ws = $wb.Worksheets.Item(77)
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 40

# Sheet 78: Morgan Park
$ws = $wb.Worksheets.Item(78)
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 47

# Sheet 79: Hyde Park
$ws = $wb.Worksheets.Item(79)
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 75

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item(84)
$ws.Range("J2").Value = 15
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 54

# Sheet 93: Mount Greenwood
$ws = $wb.Worksheets.Item(93)
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 14
